$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "323.26"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.80%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.55"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.10%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.869"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "11.37%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08021"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.92%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.85%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.628"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.26%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.947"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.34%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9293"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.80%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1271"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-6.20%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1956"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.04%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.715"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "29.09%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09102"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.44%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03593"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.42%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1049"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "9.40%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001301"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-4.03%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006269"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.13%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.44%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3540"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.60%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.39%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2447"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-4.76%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04410"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.81%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001263"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.31%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004393"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.81%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02512"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.21%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05225"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.06%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007458"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.99%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009606"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.05%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1406"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.42%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01109"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "15.45%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006746"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.71%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-9.94%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-7.76%"
